$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Consolidate the title's "A" / " " / "slide" runs into a single run.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "A slide"

# Consolidate the table cell's "a" / " " / "table" runs into a single run.
$tableShape = $s.Shapes.Item(3)
$tbl = $tableShape.Table
$cell = $tbl.Cell(1, 2)
$cell.Shape.TextFrame.TextRange.Text = "a table"
